$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Seed A15/C15 first (in left-to-right order) so the new shared-string
# entries land in the same order as the target workbook: "OpenXmlSdk" then
# the homepage link, before the row-14 formatting copy (and before the
# hyperlink's own style stamp) overwrite them.
$ws.Range("A15").Value = "OpenXmlSdk"
$ws.Range("C15").Value = "https://github.com/OfficeDev/Open-XML-SDK "
$ws.Hyperlinks.Add($ws.Range("C15"), "https://github.com/OfficeDev/Open-XML-SDK") | Out-Null

# Copy formatting from row 14 so fonts/fills/borders/wrap match the rest of
# the table (this also overwrites the stray style Hyperlinks.Add stamped on
# C15, restoring the shared "Link" style already used by column C).
$ws.Range("A14:D14").Copy($ws.Range("A15:D15"))

# Now (re)write the actual cell contents of the new row.
$ws.Range("A15").Value = "OpenXmlSdk"
$ws.Range("B15").Value = "Apache 2.0"
$ws.Range("C15").Value = "https://github.com/OfficeDev/Open-XML-SDK "
$ws.Range("D15").Value = " - Library for working with OpenXML file `n - Included as Dll, code is available on the project's homepage"

# Match row height used by the other data rows (2 lines of wrapped text).
$ws.Rows.Item(15).RowHeight = $ws.Rows.Item(14).RowHeight

# Restore the selection the author left the sheet in.
$ws.Range("B10").Select() | Out-Null
